# Generate Report for Handoff
# Adds two new "ready for handoff" entries (cddf2373-... and cf08c492-...)
# to all three sheets of the localization-status workbook: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob"

# ---------------------------------------------------------------------------
# New file identities being handed off
# ---------------------------------------------------------------------------
$file1Name = "cddf2373-82dd-4436-b709-5c1c90e92a5e.md"
$file1Path = "e2e\cddf2373-82dd-4436-b709-5c1c90e92a5e.md"
$file1Hash = "6cba01feba65b604ccd79c36d220581f6415526e"
$file1ZhCnXlf = "cddf2373-82dd-4436-b709-5c1c90e92a5e." + $file1Hash + ".zh-cn.xlf"
$file1DeDeXlf = "cddf2373-82dd-4436-b709-5c1c90e92a5e." + $file1Hash + ".de-de.xlf"

$file2Name = "cf08c492-7abe-48d9-9567-923b99d8096d.md"
$file2Path = "e2e\cf08c492-7abe-48d9-9567-923b99d8096d.md"
$file2Hash = "b0386ec8f3ed8112f4b3b40272ecd24efc48d729"
$file2ZhCnXlf = "cf08c492-7abe-48d9-9567-923b99d8096d." + $file2Hash + ".zh-cn.xlf"
$file2DeDeXlf = "cf08c492-7abe-48d9-9567-923b99d8096d." + $file2Hash + ".de-de.xlf"

$zhCnHandoffDate = "2016-08-15 16:38:57"
$deDeHandoffDate = "2016-08-15 16:39:05"
$genDate         = "2016-08-15 16:39:05"
$neverDate       = "0001-01-01 00:00:00"
$dateFormat      = "yyyy-mm-dd HH:mm:ss"

# ===========================================================================
# Sheet "Overview" -- rows 4 and 5
# ===========================================================================
$wsOverview.Cells.Item(4, 1).Value = $file1Name
$wsOverview.Cells.Item(4, 3).Value = ".md"
$wsOverview.Cells.Item(4, 4).Value = ""
$wsOverview.Cells.Item(4, 5).Value = "Ready for handoff"
$wsOverview.Cells.Item(4, 6).Value = "Ready for handoff"
$wsOverview.Cells.Item(4, 7).Value = $genDate
$wsOverview.Cells.Item(4, 7).NumberFormat = $dateFormat
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), ($repoBase + "/" + $file1Hash + "/e2e/" + $file1Name), "", "", $file1Path)

$wsOverview.Cells.Item(5, 1).Value = $file2Name
$wsOverview.Cells.Item(5, 3).Value = ".md"
$wsOverview.Cells.Item(5, 4).Value = ""
$wsOverview.Cells.Item(5, 5).Value = "Ready for handoff"
$wsOverview.Cells.Item(5, 6).Value = "Ready for handoff"
$wsOverview.Cells.Item(5, 7).Value = $genDate
$wsOverview.Cells.Item(5, 7).NumberFormat = $dateFormat
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), ($repoBase + "/" + $file2Hash + "/e2e/" + $file2Name), "", "", $file2Path)

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G5"))

# ===========================================================================
# Sheet "zh-cn" -- rows 4 and 5
# ===========================================================================
$wsZhCn.Cells.Item(4, 2).Value  = ".md"
$wsZhCn.Cells.Item(4, 3).Value  = "Ready for handoff"
$wsZhCn.Cells.Item(4, 4).Value  = "e2e"
$wsZhCn.Cells.Item(4, 5).Value  = "ht"
$wsZhCn.Cells.Item(4, 6).Value  = "False"
$wsZhCn.Cells.Item(4, 7).Value  = $file1ZhCnXlf
$wsZhCn.Cells.Item(4, 8).Value  = $zhCnHandoffDate
$wsZhCn.Cells.Item(4, 8).NumberFormat = $dateFormat
$wsZhCn.Cells.Item(4, 9).Value  = ""
$wsZhCn.Cells.Item(4, 10).Value = ""
$wsZhCn.Cells.Item(4, 11).Value = $neverDate
$wsZhCn.Cells.Item(4, 11).NumberFormat = $dateFormat
$wsZhCn.Cells.Item(4, 12).Value = ""
$wsZhCn.Cells.Item(4, 13).Value = "True"
$wsZhCn.Cells.Item(4, 14).Value = ""
$wsZhCn.Cells.Item(4, 15).Value = "False"
$wsZhCn.Cells.Item(4, 16).Value = ""
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), ($repoBase + "/" + $file1Hash + "/e2e/" + $file1Name), "", "", $file1Name)

$wsZhCn.Cells.Item(5, 2).Value  = ".md"
$wsZhCn.Cells.Item(5, 3).Value  = "Ready for handoff"
$wsZhCn.Cells.Item(5, 4).Value  = "e2e"
$wsZhCn.Cells.Item(5, 5).Value  = "ht"
$wsZhCn.Cells.Item(5, 6).Value  = "False"
$wsZhCn.Cells.Item(5, 7).Value  = $file2ZhCnXlf
$wsZhCn.Cells.Item(5, 8).Value  = $zhCnHandoffDate
$wsZhCn.Cells.Item(5, 8).NumberFormat = $dateFormat
$wsZhCn.Cells.Item(5, 9).Value  = ""
$wsZhCn.Cells.Item(5, 10).Value = ""
$wsZhCn.Cells.Item(5, 11).Value = $neverDate
$wsZhCn.Cells.Item(5, 11).NumberFormat = $dateFormat
$wsZhCn.Cells.Item(5, 12).Value = ""
$wsZhCn.Cells.Item(5, 13).Value = "True"
$wsZhCn.Cells.Item(5, 14).Value = ""
$wsZhCn.Cells.Item(5, 15).Value = "False"
$wsZhCn.Cells.Item(5, 16).Value = ""
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), ($repoBase + "/" + $file2Hash + "/e2e/" + $file2Name), "", "", $file2Name)

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P5"))

# ===========================================================================
# Sheet "de-de" -- rows 4 and 5
# ===========================================================================
$wsDeDe.Cells.Item(4, 2).Value  = ".md"
$wsDeDe.Cells.Item(4, 3).Value  = "Ready for handoff"
$wsDeDe.Cells.Item(4, 4).Value  = "e2e"
$wsDeDe.Cells.Item(4, 5).Value  = "ht"
$wsDeDe.Cells.Item(4, 6).Value  = "False"
$wsDeDe.Cells.Item(4, 7).Value  = $file1DeDeXlf
$wsDeDe.Cells.Item(4, 8).Value  = $deDeHandoffDate
$wsDeDe.Cells.Item(4, 8).NumberFormat = $dateFormat
$wsDeDe.Cells.Item(4, 9).Value  = ""
$wsDeDe.Cells.Item(4, 10).Value = ""
$wsDeDe.Cells.Item(4, 11).Value = $neverDate
$wsDeDe.Cells.Item(4, 11).NumberFormat = $dateFormat
$wsDeDe.Cells.Item(4, 12).Value = ""
$wsDeDe.Cells.Item(4, 13).Value = "True"
$wsDeDe.Cells.Item(4, 14).Value = ""
$wsDeDe.Cells.Item(4, 15).Value = "False"
$wsDeDe.Cells.Item(4, 16).Value = ""
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), ($repoBase + "/" + $file1Hash + "/e2e/" + $file1Name), "", "", $file1Name)

$wsDeDe.Cells.Item(5, 2).Value  = ".md"
$wsDeDe.Cells.Item(5, 3).Value  = "Ready for handoff"
$wsDeDe.Cells.Item(5, 4).Value  = "e2e"
$wsDeDe.Cells.Item(5, 5).Value  = "ht"
$wsDeDe.Cells.Item(5, 6).Value  = "False"
$wsDeDe.Cells.Item(5, 7).Value  = $file2DeDeXlf
$wsDeDe.Cells.Item(5, 8).Value  = $deDeHandoffDate
$wsDeDe.Cells.Item(5, 8).NumberFormat = $dateFormat
$wsDeDe.Cells.Item(5, 9).Value  = ""
$wsDeDe.Cells.Item(5, 10).Value = ""
$wsDeDe.Cells.Item(5, 11).Value = $neverDate
$wsDeDe.Cells.Item(5, 11).NumberFormat = $dateFormat
$wsDeDe.Cells.Item(5, 12).Value = ""
$wsDeDe.Cells.Item(5, 13).Value = "True"
$wsDeDe.Cells.Item(5, 14).Value = ""
$wsDeDe.Cells.Item(5, 15).Value = "False"
$wsDeDe.Cells.Item(5, 16).Value = ""
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), ($repoBase + "/" + $file2Hash + "/e2e/" + $file2Name), "", "", $file2Name)

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P5"))

Write-Host "Report generated for handoff: added 2 rows to Overview, zh-cn, de-de"
